$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 3526.3333
$ws.Range("I42").Value = 4030
$ws.Range("J42").Value = 3274.5
$ws.Range("K42").Value = 12090
$ws.Range("L42").Value = 9823.5
$ws.Range("M42").Value = -11860
$ws.Range("N42").Value = -10283.5
$ws.Range("H43").Value = 2212.1428
$ws.Range("I43").Value = 1029.3334
$ws.Range("J43").Value = 3099.25
$ws.Range("K43").Value = 1029.3334
$ws.Range("L43").Value = 3099.25
$ws.Range("M43").Value = -960.3334
$ws.Range("N43").Value = -3237.25
$ws.Range("H111").Value = 1445.7059
$ws.Range("J111").Value = 2438.5
$ws.Range("L111").Value = 7315.5
$ws.Range("N111").Value = -13449.5
$ws.Range("H112").Value = 839082.1
$ws.Range("I112").Value = 3192
$ws.Range("J112").Value = 991062.2
$ws.Range("K112").Value = 9576
$ws.Range("L112").Value = 2973186.6
$ws.Range("M112").Value = -8468
$ws.Range("N112").Value = -2975402.6
$ws.Range("H127").Value = 1597
$ws.Range("I127").Value = 1597
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 4791
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 169
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 4994
$ws.Range("J19").Value = 4994
$ws.Range("L19").Value = 4994
$ws.Range("N19").Value = -5452
$ws.Range("H22").Value = 1417.3636
$ws.Range("J22").Value = 1980
$ws.Range("L22").Value = 1980
$ws.Range("N22").Value = -2578
$ws.Range("H41").Value = 17391
$ws.Range("I41").Value = 1449
$ws.Range("K41").Value = 1449
$ws.Range("M41").Value = -1035
$ws.Range("H45").Value = 2190.3
$ws.Range("I45").Value = 836.1
$ws.Range("K45").Value = 836.1
$ws.Range("M45").Value = -459.1
$ws.Range("H74").Value = 2229.8518
$ws.Range("I74").Value = 2666.125
$ws.Range("J74").Value = 1595.2727
$ws.Range("K74").Value = 2666.125
$ws.Range("L74").Value = 1595.2727
$ws.Range("M74").Value = -1792.125
$ws.Range("N74").Value = -3343.2727
$ws.Range("H77").Value = 2229.8518
$ws.Range("I77").Value = 2666.125
$ws.Range("J77").Value = 1595.2727
$ws.Range("K77").Value = 13330.625
$ws.Range("L77").Value = 7976.363499999999
$ws.Range("M77").Value = -8962.625
$ws.Range("N77").Value = -16712.3635
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1276.5
$ws.Range("I8").Value = 869
$ws.Range("K8").Value = 869
$ws.Range("M8").Value = -729
$ws.Range("H25").Value = 8310.429
$ws.Range("I25").Value = 7014
$ws.Range("K25").Value = 7014
$ws.Range("M25").Value = -6779
$ws.Range("H107").Value = 5206.4287
$ws.Range("I107").Value = 3505.5
$ws.Range("K107").Value = 3505.5
$ws.Range("M107").Value = -1585.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 823.2143
$ws.Range("I16").Value = 728.4167
$ws.Range("J16").Value = 1392
$ws.Range("K16").Value = 728.4167
$ws.Range("L16").Value = 1392
$ws.Range("M16").Value = -441.4167
$ws.Range("N16").Value = -1966
$ws.Range("H18").Value = 21786
$ws.Range("J18").Value = 21786
$ws.Range("L18").Value = 21786
$ws.Range("N18").Value = -22246
$ws.Range("H107").Value = 3552.2666
$ws.Range("J107").Value = 3600.8333
$ws.Range("L107").Value = 3600.8333
$ws.Range("N107").Value = -7440.8333
$ws.Range("H113").Value = 823.2143
$ws.Range("I113").Value = 728.4167
$ws.Range("J113").Value = 1392
$ws.Range("K113").Value = 728.4167
$ws.Range("L113").Value = 1392
$ws.Range("M113").Value = 1441.5833
$ws.Range("N113").Value = -5732
$ws.Range("H119").Value = 64993.5
$ws.Range("J119").Value = 64993.5
$ws.Range("L119").Value = 64993.5
$ws.Range("N119").Value = -74669.5
$ws.Range("H134").Value = 1780.0322
$ws.Range("I134").Value = 1613.6786
$ws.Range("K134").Value = 4841.0358
$ws.Range("M134").Value = -2306.0358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 350.05554
$ws.Range("J40").Value = 377.36365
$ws.Range("L40").Value = 1509.4546
$ws.Range("N40").Value = -1647.4546
$ws.Range("H42").Value = 10498.5
$ws.Range("I42").Value = 8998
$ws.Range("J42").Value = 15000
$ws.Range("K42").Value = 26994
$ws.Range("L42").Value = 45000
$ws.Range("M42").Value = -26460
$ws.Range("N42").Value = -46068
$ws.Range("H61").Value = 2159.6667
$ws.Range("I61").Value = 2688.2856
$ws.Range("J61").Value = 309.5
$ws.Range("K61").Value = 8064.8568
$ws.Range("L61").Value = 928.5
$ws.Range("M61").Value = -7849.8568
$ws.Range("N61").Value = -1358.5
$ws.Range("H109").Value = 5806.3335
$ws.Range("I109").Value = 5806.3335
$ws.Range("K109").Value = 17419.0005
$ws.Range("M109").Value = -16379.0005
$ws.Range("H129").Value = 2245.2727
$ws.Range("I129").Value = 325.66666
$ws.Range("J129").Value = 2965.125
$ws.Range("K129").Value = 976.9999799999999
$ws.Range("L129").Value = 8895.375
$ws.Range("M129").Value = 4023.00002
$ws.Range("N129").Value = -18895.375
$ws.Range("H131").Value = 1502.4667
$ws.Range("J131").Value = 1513.1389
$ws.Range("L131").Value = 4539.4167
$ws.Range("N131").Value = -14619.4167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 9834.75
$ws.Range("J15").Value = 9834.75
$ws.Range("L15").Value = 9834.75
$ws.Range("N15").Value = -10410.75
$ws.Range("H70").Value = 12937.3125
$ws.Range("I70").Value = 11999.917
$ws.Range("J70").Value = 15749.5
$ws.Range("K70").Value = 11999.917
$ws.Range("L70").Value = 15749.5
$ws.Range("M70").Value = -11729.917
$ws.Range("N70").Value = -16289.5
$ws.Range("H73").Value = 12937.3125
$ws.Range("I73").Value = 11999.917
$ws.Range("J73").Value = 15749.5
$ws.Range("K73").Value = 11999.917
$ws.Range("L73").Value = 15749.5
$ws.Range("M73").Value = -11063.917
$ws.Range("N73").Value = -17621.5
$ws.Range("H81").Value = 9834.75
$ws.Range("J81").Value = 9834.75
$ws.Range("L81").Value = 9834.75
$ws.Range("N81").Value = -11830.75
$ws.Range("H84").Value = 9834.75
$ws.Range("J84").Value = 9834.75
$ws.Range("L84").Value = 29504.25
$ws.Range("N84").Value = -39488.25
$ws.Range("H113").Value = 3324.75
$ws.Range("I113").Value = 3033
$ws.Range("J113").Value = 3499.8
$ws.Range("K113").Value = 3033
$ws.Range("L113").Value = 3499.8
$ws.Range("M113").Value = -863
$ws.Range("N113").Value = -7839.8
$ws.Range("H132").Value = 3408
$ws.Range("I132").Value = 4724.5
$ws.Range("J132").Value = 2749.75
$ws.Range("K132").Value = 14173.5
$ws.Range("L132").Value = 8249.25
$ws.Range("M132").Value = -11643.5
$ws.Range("N132").Value = -13309.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 8032
$ws.Range("I9").Value = 2052
$ws.Range("J9").Value = 19992
$ws.Range("K9").Value = 2052
$ws.Range("L9").Value = 19992
$ws.Range("M9").Value = -1828
$ws.Range("N9").Value = -20440
$ws.Range("H10").Value = 8197.333000000001
$ws.Range("J10").Value = 9761.333000000001
$ws.Range("L10").Value = 9761.333000000001
$ws.Range("N10").Value = -10041.333
$ws.Range("H12").Value = 14331.444
$ws.Range("J12").Value = 13623.25
$ws.Range("L12").Value = 13623.25
$ws.Range("N12").Value = -13963.25
$ws.Range("H34").Value = 24333.334
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H46").Value = 2956.6667
$ws.Range("I46").Value = 632.6667
$ws.Range("K46").Value = 632.6667
$ws.Range("M46").Value = -444.6667
$ws.Range("H75").Value = 49998.5
$ws.Range("J75").Value = 49998.5
$ws.Range("L75").Value = 49998.5
$ws.Range("N75").Value = -51870.5
$ws.Range("H76").Value = 9997
$ws.Range("J76").Value = 9997
$ws.Range("L76").Value = 9997
$ws.Range("N76").Value = -10673
$ws.Range("H78").Value = 49998.5
$ws.Range("J78").Value = 49998.5
$ws.Range("L78").Value = 149995.5
$ws.Range("N78").Value = -159355.5
$ws.Range("H79").Value = 9997
$ws.Range("J79").Value = 9997
$ws.Range("L79").Value = 9997
$ws.Range("N79").Value = -12337
$ws.Range("H103").Value = 12600.2
$ws.Range("J103").Value = 12600.2
$ws.Range("L103").Value = 12600.2
$ws.Range("N103").Value = -14944.2
$ws.Range("H132").Value = 8175.75
$ws.Range("J132").Value = 4602
$ws.Range("L132").Value = 13806
$ws.Range("N132").Value = -18866

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 346786.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 346786.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 346786.25
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -347566.25
$ws.Range("H46").Value = 94952
$ws.Range("J46").Value = 94952
$ws.Range("L46").Value = 94952
$ws.Range("N46").Value = -95414
$ws.Range("H134").Value = 94952
$ws.Range("J134").Value = 94952
$ws.Range("L134").Value = 284856
$ws.Range("N134").Value = -289926
$ws.Range("H136").Value = 743.35486
$ws.Range("I136").Value = 653.3103599999999
$ws.Range("K136").Value = 1959.93108
$ws.Range("M136").Value = 590.0689200000002
